$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Manager" -> "Manager id" (project managers are now tracked by id)
$ws.Range("K1").Value = "Manager id"

# Data row: replace the manager's name ("Jessica") with their numeric id
$ws.Range("K2").Value = 9

# Match the reviewer's on-screen state when the check was made: zoomed in,
# scrolled over to the manager/officer columns, with K2 selected.
$excel.ActiveWindow.Zoom = 156
$ws.Range("K2").Select()
